$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.497.68"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.59%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.747.22"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.50%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.73%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.19"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.90%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.64%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4441"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +4.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3594"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.62%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07496"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.17%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.96"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -5.16%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.092"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.43%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.93%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.72"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -4.00%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.015"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.02%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.100"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.93%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.751.05"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.07%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.25"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.84%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001060"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.03%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06386"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.80%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.57%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.79"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.31%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.817"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.20%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.558.22"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.35%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.18"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.58%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.097"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.61%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.30"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.70%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.50"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.949.45"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.01%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.073"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -4.37%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.54"

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.078"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -7.54%  "

$ws.Range("B32").Value = "HuobiToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.666"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +4.62%  "

$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09054"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.16%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.532"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.71%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.93"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -5.51%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02287"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.27%  "

$ws.Range("B37").Value = "Algorand"
$ws.Range("C37").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2089"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.03%  "

$ws.Range("B38").Value = "TheSandbox"
$ws.Range("C38").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.6353"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.23%  "

$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05996"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.89%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.932"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.54%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.202"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.12%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.385"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.19%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.770"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.12%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.30"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.43%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.719"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.90%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5868"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.66%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "121.84"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.83%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.949"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.54%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.143"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.60%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06849"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.87%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.19"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.11%  "
